$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/link/percentage updates (not numeric-looking, safe as-is)
$ws.Range("D2").Value = '90.780.85'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '3.148.12'
$ws.Range("E3").Value = '  +2.60%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("E7").Value = '  +29.86%  '
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.147.12'
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("E11").Value = '  +11.52%  '
$ws.Range("E12").Value = '  +7.32%  '
$ws.Range("E13").Value = '  +6.52%  '
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("D16").Value = '90.663.42'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '3.743.76'
$ws.Range("E17").Value = '  +3.06%  '
$ws.Range("D18").Value = '3.134.10'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("E19").Value = '  +8.84%  '
$ws.Range("E20").Value = '  +6.51%  '
$ws.Range("E21").Value = '  +10.05%  '
$ws.Range("E22").Value = '  -6.66%  '
$ws.Range("E23").Value = '  +8.28%  '
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("E25").Value = '  +14.76%  '
$ws.Range("E26").Value = '  +4.45%  '
$ws.Range("E27").Value = '  +3.74%  '
$ws.Range("D28").Value = '3.332.12'
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -3.42%  '
$ws.Range("E31").Value = '  +5.48%  '
$ws.Range("E32").Value = '  -5.96%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E33").Value = '  +19.08%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E34").Value = '  +45.75%  '
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E36").Value = '  +5.56%  '
$ws.Range("E37").Value = '  +5.58%  '
$ws.Range("E38").Value = '  -6.62%  '
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("E41").Value = '  +27.36%  '
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("E44").Value = '  +14.05%  '
$ws.Range("E45").Value = '  +5.91%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E47").Value = '  +19.24%  '
$ws.Range("E48").Value = '  +4.21%  '
$ws.Range("E49").Value = '  +8.01%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E50").Value = '  +4.22%  '
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E51").Value = '  +10.04%  '

# Price values that look numeric must be forced to remain text,
# matching the original inlineStr/shared-string text cells, and then
# the style is reset to Normal so no stray number-format style is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.16'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.752'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.201'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '474.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000210'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.202'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '516.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.144'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0911'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.418'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.717'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '150.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("D51").Style = "Normal"
